# update list dir logic
#
# Append a new log row for the latest directory-listing check run. The new
# row duplicates the most recent prior run (row 5) column-for-column -- using
# a range Copy so every value's original type (text vs. number) is preserved
# verbatim, including date-like text in column B that a plain .Value write
# would otherwise have Excel auto-convert into a real date serial -- and then
# stamps it with this run's own "Checking Time".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRange = $ws.Range("A5:BF5")
$dstRange = $ws.Range("A6:BF6")
$srcRange.Copy($dstRange)

# The new run's "Checking Time" differs from the row it was copied from.
$ws.Cells.Item(6, 3).Value = "10:33:28"
